$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.08824
$ws.Range("C2").Value = 1.5933402
$ws.Range("E2").Value = 0.06403226187196739

$ws.Range("B3").Value = 0.15038
$ws.Range("C3").Value = 1.4004756
$ws.Range("E3").Value = 0.06403226187196739

$ws.Range("B4").Value = 0.21124
$ws.Range("C4").Value = 1.2043737
$ws.Range("E4").Value = 0.06403226187196739
